$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section
# ("LOQ4205: Sistemas Produtivos II (Requisito fraco)").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4205*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The three paragraphs immediately following it are the footer block that
    # was scraped from the site (a blank line, the "Ver no Jupiter..." line,
    # and the copyright/Jekyll line). Remove all three, leaving the blank
    # paragraph (and page break paragraph) that originally followed them.
    $first = $target.Next()
    $second = $first.Next()
    $third = $second.Next()

    $deleteRange = $d.Range($first.Range.Start, $third.Range.End)
    $deleteRange.Delete()
}
